$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.253.11'
$ws.Range('E2').Value = '  +6.32%  '
$ws.Range('D3').Value = '2.593.28'
$ws.Range('E3').Value = '  +6.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '182.67'
$ws.Range('E5').Value = '  +12.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '578.97'
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('E8').Value = '  +3.43%  '
$ws.Range('E9').Value = '  +14.73%  '
$ws.Range('D10').Value = '2.593.00'
$ws.Range('E10').Value = '  +6.50%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.358'
$ws.Range('E12').Value = '  +7.76%  '
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = '73.195.42'
$ws.Range('E14').Value = '  +6.44%  '
$ws.Range('D15').Value = '3.067.72'
$ws.Range('E15').Value = '  +6.40%  '
$ws.Range('E16').Value = '  +4.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.77'
$ws.Range('E17').Value = '  +11.07%  '
$ws.Range('D18').Value = '2.586.88'
$ws.Range('E18').Value = '  +6.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.80'
$ws.Range('E19').Value = '  +27.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.72'
$ws.Range('E20').Value = '  +11.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '368.30'
$ws.Range('E21').Value = '  +8.64%  '
$ws.Range('E22').Value = '  +14.02%  '
$ws.Range('E23').Value = '  +5.68%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.19'
$ws.Range('E25').Value = '  +3.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.09'
$ws.Range('E26').Value = '  +10.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.16'
$ws.Range('E27').Value = '  +11.96%  '
$ws.Range('D28').Value = '2.713.36'
$ws.Range('E28').Value = '  +5.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +12.13%  '
$ws.Range('E31').Value = '  +17.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '500.20'
$ws.Range('E32').Value = '  +16.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.58'
$ws.Range('E33').Value = '  +6.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.72'
$ws.Range('E34').Value = '  +6.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +12.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.23'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.97'
$ws.Range('E38').Value = '  +5.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.23'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.81'
$ws.Range('E41').Value = '  +10.84%  '
$ws.Range('E42').Value = '  +9.50%  '
$ws.Range('E43').Value = '  +7.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '154.74'
$ws.Range('E44').Value = '  +19.39%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0860'
$ws.Range('E45').Value = '  +19.52%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.16'
$ws.Range('E46').Value = '  +7.92%  '
$ws.Range('E47').Value = '  +13.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '38.41'
$ws.Range('E48').Value = '  +2.64%  '
$ws.Range('E49').Value = '  +7.36%  '
$ws.Range('E50').Value = '  +8.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.07'
$ws.Range('E51').Value = '  +19.09%  '
